$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells being updated to Text format so that
# values such as "135.77", "18.00" or "8.10" are preserved exactly as
# text (matching the source data) instead of being parsed as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D9", "D14", "D15", "D16", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D29", "D33", "D35", "D38", "D40", "D41", "D42", "D44", "D45", "D46", "D50")
foreach ($pc in $priceCells) {
    $ws.Range($pc).NumberFormat = "@"
}

$ws.Range("D2").Value = "57.011.42"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "2.341.55"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "518.38"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "135.77"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "2.353.32"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").Value = "23.95"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "2.758.87"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").Value = "56.991.63"
$ws.Range("E16").Value = "  +1.32%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "2.350.56"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").Value = "  +0.74%  "
$ws.Range("D20").Value = "326.80"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "6.75"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "61.27"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "0.165"
$ws.Range("E25").Value = "  +4.97%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").Value = "8.10"
$ws.Range("E27").Value = "  +6.10%  "
$ws.Range("E28").Value = "  +9.34%  "
$ws.Range("D29").Value = "170.14"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "18.57"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("E39").Value = "  +3.61%  "
$ws.Range("D40").Value = "38.47"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("D41").Value = "149.19"
$ws.Range("E41").Value = "  +6.42%  "
$ws.Range("D42").Value = "0.383"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("D44").Value = "280.45"
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("D45").Value = "5.20"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "0.0936"
$ws.Range("E46").Value = "  +1.12%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +1.98%  "
$ws.Range("D50").Value = "18.00"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("E51").Value = "  -0.08%  "
